$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 13 - this shifts the existing row 13
# ("Телевизор", 2) down to row 14 (it keeps its original, unstyled
# formatting), and the new blank row 13 inherits the formatting of the
# row above it (row 12, which uses the "data row" style shared by A:B).
$ws.Rows("13:13").Insert()

# Restore the original row's data into the now-styled row 13.
$ws.Range("A13").Value = "Телевизор"
$ws.Range("B13").Value = 2

# Write the teacher's new data row into row 14 (shifted-down row, which
# kept the original - unstyled - formatting).
$ws.Range("A14").Value = "Стиральная машина"
$ws.Range("B14").Value = 1.39
